# Applies the coin-collection "exchange status" updates described by the
# commit: a batch of cells across several sheets are flipped from 0 -> 1
# (or 1 -> 2 where a coin can now be exchanged twice), plus the sheet that
# is on top / selected in the workbook changes from
# "The_Beautiful_Quarters" to "American_Innovation_Dollars".

$wb = $excel.ActiveWorkbook

# Worksheet index map (1-based, matches xl/workbook.xml sheet order):
#  1 "Linkoln_Cents"
#  2 "Westward_Journey_Nickels"
#  3 "State&Territorial_Quarters"
#  4 "The_Beautiful_Quarters"
#  5 "American_Women_Quarters"
#  6 "Presidential_Dollars"
#  7 "Sacagawea&Native_Dollars"
#  8 "American_Innovation_Dollars"
#  9 "Links"

$wsBeautifulQuarters   = $wb.Worksheets.Item(4)
$wsAmericanWomen       = $wb.Worksheets.Item(5)
$wsPresidentialDollars = $wb.Worksheets.Item(6)
$wsSacagawea           = $wb.Worksheets.Item(7)
$wsAmericanInnovation  = $wb.Worksheets.Item(8)

# --- "The_Beautiful_Quarters" -------------------------------------------
$wsBeautifulQuarters.Range("J33").Value = 1
$wsBeautifulQuarters.Range("I34").Value = 1
$wsBeautifulQuarters.Range("J34").Value = 1
$wsBeautifulQuarters.Range("J40").Value = 1
$wsBeautifulQuarters.Range("K56").Value = 1
$wsBeautifulQuarters.Range("I57").Value = 1
$wsBeautifulQuarters.Range("J58").Value = 1

# --- "American_Women_Quarters" ------------------------------------------
$wsAmericanWomen.Range("H6").Value = 1
$wsAmericanWomen.Range("H7").Value = 1
$wsAmericanWomen.Range("H8").Value = 1
$wsAmericanWomen.Range("H9").Value = 1
$wsAmericanWomen.Range("H12").Value = 1
$wsAmericanWomen.Range("I12").Value = 1
$wsAmericanWomen.Range("H13").Value = 1
$wsAmericanWomen.Range("I13").Value = 1

# --- "Presidential_Dollars" ----------------------------------------------
$wsPresidentialDollars.Range("G26").Value = 1
$wsPresidentialDollars.Range("G42").Value = 1
$wsPresidentialDollars.Range("H42").Value = 1

# --- "Sacagawea&Native_Dollars" -------------------------------------------
$wsSacagawea.Range("G23").Value = 2
$wsSacagawea.Range("G24").Value = 1
$wsSacagawea.Range("G26").Value = 1

# --- "American_Innovation_Dollars" ----------------------------------------
$wsAmericanInnovation.Range("G13").Value = 1
$wsAmericanInnovation.Range("H13").Value = 1
$wsAmericanInnovation.Range("G14").Value = 1
$wsAmericanInnovation.Range("H14").Value = 1
$wsAmericanInnovation.Range("H15").Value = 1
$wsAmericanInnovation.Range("H18").Value = 1
$wsAmericanInnovation.Range("H19").Value = 1
$wsAmericanInnovation.Range("G20").Value = 2
$wsAmericanInnovation.Range("G21").Value = 1

# --- Update the remembered cursor/selection on each touched sheet --------
$wsBeautifulQuarters.Range("J59").Select() | Out-Null
$wsAmericanWomen.Range("J14").Select() | Out-Null
$wsPresidentialDollars.Range("F38").Select() | Out-Null
$wsSacagawea.Range("I29").Select() | Out-Null
$wsAmericanInnovation.Range("I18").Select() | Out-Null

# The workbook now opens on "American_Innovation_Dollars" (tab index 7,
# 0-based) instead of "The_Beautiful_Quarters" (tab index 3, 0-based).
$wsAmericanInnovation.Activate()
